# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Mon Oct  9 11:34:44 UTC 2023 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D:E data range to Text format before writing, so that numeric-looking
# strings (e.g. "207.65", "0.0968") are preserved as text instead of being
# auto-converted to numbers by Excel -- matching the original inlineStr cells.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.559.12"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.595.45"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "207.65"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "0.500"
$ws.Range("E6").Value = "  -4.03%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "22.24"
$ws.Range("E8").Value = "  -4.22%  "
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "0.0588"
$ws.Range("E10").Value = "  -3.49%  "
$ws.Range("D11").Value = "0.0871"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "1.820.91"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "1.592.40"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("D16").Value = "63.39"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "27.516.48"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "216.79"
$ws.Range("E18").Value = "  -5.13%  "
$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").Value = "9.71"
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").Value = "154.65"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -2.30%  "
$ws.Range("D28").Value = "15.01"
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("E29").Value = "  -4.49%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "1.349.15"
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "0.960"
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("D39").Value = "0.538"
$ws.Range("E39").Value = "  -2.77%  "
$ws.Range("D40").Value = "0.815"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "0.962"
$ws.Range("E42").Value = "  -5.08%  "
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").Value = "63.86"
$ws.Range("E44").Value = "  -2.55%  "
$ws.Range("E45").Value = "  -4.23%  "
$ws.Range("D46").Value = "1.732.00"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "87.17"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.0968"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0496"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0958"
$ws.Range("E51").Value = "  -6.60%  "

# Remove the temporary text-format override so cells end up with the default
# (unstyled) formatting again, same as the rest of the sheet.
$dataRange.ClearFormats()

